# Auto-generated from the canonical OOXML diff.
# Updates the "Price" (D) and "Volume(1h)" (E) columns for the cryptos table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.943.57"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.638.16"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'212.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'23.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'0.0882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.869.88"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.639.45"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'0.568"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "27.948.26"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'231.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").Value = "'153.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "'6.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'15.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'0.0483"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "1.401.83"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D38").Value = "'0.563"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'0.927"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "'0.876"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'66.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "1.779.25"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'87.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "'0.100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "'0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'7.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.08%  "
